$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '76.590.02'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  +0.90%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.027.68'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +4.20%  '

# Row 4
$ws.Range("E4").Value = '  -0.01%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '201.67'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +1.13%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '629.81'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +5.41%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +0.05%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.555'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +1.41%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.212'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +7.04%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '3.027.28'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +4.27%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.437'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +2.54%  '

# Row 12
$ws.Range("E12").Value = '  -0.27%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.18'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +6.86%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.585.59'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +4.28%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '29.45'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +7.29%  '

# Row 16
$ws.Range("B16").Value = 'ShibaInu'
$ws.Range("C16").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000195'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +2.63%  '

# Row 17
$ws.Range("B17").Value = 'WrappedBTC'
$ws.Range("C17").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '76.459.53'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +0.88%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.021.17'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +4.17%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.52'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +6.23%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '9.00'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +1.11%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '376.29'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +1.13%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.37'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +2.00%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.29'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +0.21%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '73.33'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +3.38%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.185.23'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +4.12%  '

# Row 26
$ws.Range("B26").Value = 'NEARProtocol'
$ws.Range("C26").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '4.41'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +5.58%  '

# Row 27
$ws.Range("B27").Value = 'Dai'
$ws.Range("C27").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.999'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +0.02%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.97'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +3.53%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.0000111'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +3.93%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.996'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -0.20%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.36'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +8.72%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.43'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +2.01%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '508.22'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +1.47%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.96'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +8.13%  '

# Row 35
$ws.Range("E35").Value = '  +0.05%  '

# Row 36
$ws.Range("E36").Value = '  +2.98%  '

# Row 37
$ws.Range("E37").Value = '  -0.95%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.385'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +11.71%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '20.02'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +1.94%  '

# Row 40
$ws.Range("B40").Value = 'Cronos'
$ws.Range("C40").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.106'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +5.67%  '

# Row 41
$ws.Range("B41").Value = 'Aave'
$ws.Range("C41").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '187.88'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +3.95%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.113'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +0.56%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.17'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +4.42%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '42.46'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +5.63%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.69'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +2.69%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.26'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +6.46%  '

# Row 48
$ws.Range("B48").Value = 'dogwifhat'
$ws.Range("C48").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.43'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +4.89%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.607'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +6.28%  '

# Row 50
$ws.Range("B50").Value = 'Mantle'
$ws.Range("C50").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.715'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +9.31%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '3.91'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +5.73%  '
